$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the two new header cells, copying the formatting used by the
# existing header row (bold, centered, bordered - style of H1).
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Fill in the new numeric columns I and J for rows 2-26.
$values = @(
    @(6, 7),
    @(7, 7),
    @(7, 7),
    @(7, 7),
    @(5, 9),
    @(6, 8),
    @(8, 8),
    @(4, 5),
    @(6, 7),
    @(8, 8),
    @(8, 8),
    @(8, 9),
    @(8, 8),
    @(3, 5),
    @(10, 10),
    @(6, 8),
    @(7, 7),
    @(5, 8),
    @(8, 9),
    @(7, 9),
    @(8, 9),
    @(8, 9),
    @(6, 7),
    @(1, 2),
    @(6, 6)
)

for ($idx = 0; $idx -lt $values.Length; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $values[$idx][0]
    $ws.Cells.Item($row, 10).Value = $values[$idx][1]
}
